$d = $word.ActiveDocument

# The document ends with a run of trailing placeholder paragraphs: a
# paragraph that only holds the auto-generated "_GoBack" bookmark, an
# empty bold/size-30 paragraph, and two completely empty paragraphs.
# None of them carry visible text or inline content. Find the last
# paragraph that actually has content (text beyond the paragraph mark,
# or an inline picture) and drop everything the document has after it,
# so the body flows straight from that paragraph into the section break.

$count = $d.Paragraphs.Count

$lastContentIndex = 0
for ($i = $count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    $r = $p.Range
    $hasShape = $r.InlineShapes.Count -gt 0
    $hasText = $r.Text.Trim().Length -gt 0
    if ($hasShape -or $hasText) {
        $lastContentIndex = $i
        break
    }
}

if ($lastContentIndex -gt 0 -and $lastContentIndex -lt $count) {
    $firstTrailing = $d.Paragraphs.Item($lastContentIndex + 1)
    $lastTrailing = $d.Paragraphs.Item($count)
    $delRange = $d.Range($firstTrailing.Range.Start, $lastTrailing.Range.End)
    $delRange.Delete()
}
